# "Generate Report for Archive"
#
# Updates the localization-status report:
#   - Status text "Ready for handoff" -> "In Translation" on all sheets
#   - Narrow the now-shorter "Status"/language columns to fit the new text

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New column width (raw OOXML character-width units) used for the
# status columns after the text got shorter. ColumnWidth on this engine
# is quantized to the nearest pixel (1/6 of a character unit), so we pick
# the ColumnWidth input that lands closest to the target raw width.
$targetRawWidth = 13.4101845877511
$newColumnWidth = ($targetRawWidth - 5/6)

# --- Overview sheet: columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedOverview = $wsOverview.UsedRange
for ($r = 1; $r -le $usedOverview.Rows.Count; $r++) {
    foreach ($colIdx in 5,6) {
        $cell = $wsOverview.Cells.Item($r, $colIdx)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn and de-de sheets: column C (Status) ---
foreach ($sheetName in "zh-cn","de-de") {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $newColumnWidth
}

$wb.Save()
Write-Host "Updated status text and column widths on Overview, zh-cn, de-de."
